$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1046
$ws1.Range("F3").Value = 19

# Sheet "全部类型"
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 1046
$ws2.Range("F3").Value = 19
